# Update workbook for "Add data for 2022-07-06" commit:
#  - Rename sheet / update header label from "through June 27" to "through June 28"
#  - Update the carjacking counts for several neighborhood/month cells

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Sheet tab name: 2022-06-27 -> 2022-06-28
$ws.Name = "Through 2022-06-28"

# Column header label (shared string) for the "June 2022" column
$ws.Range("B1").Value = "June 2022 (through June 28)"

# Updated / new counts
$ws.Range("B2").Value = 9
$ws.Range("H3").Value = 3
$ws.Range("N3").Value = 4
$ws.Range("Z3").Value = 3
$ws.Range("N4").Value = 4
$ws.Range("H6").Value = 6
$ws.Range("N6").Value = 4
$ws.Range("T6").Value = 2
$ws.Range("AF6").Value = 3
$ws.Range("Z9").Value = 5
$ws.Range("H10").Value = 6
$ws.Range("N10").Value = 7
$ws.Range("AL10").Value = 5
$ws.Range("B12").Value = 6
$ws.Range("Z13").Value = 1
$ws.Range("AF14").Value = 4
$ws.Range("N18").Value = 1
$ws.Range("B20").Value = 4
$ws.Range("AF24").Value = 1
$ws.Range("AL24").Value = 2
$ws.Range("T30").Value = 1
$ws.Range("N39").Value = 2
$ws.Range("N47").Value = 2
$ws.Range("AF50").Value = 1
$ws.Range("T51").Value = 2
$ws.Range("AL59").Value = 1
$ws.Range("N70").Value = 1
$ws.Range("AL92").Value = 2
$ws.Range("AF95").Value = 2
